$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Rectification calcul rotation: update the two input angles and the
# wheel/robot diameter so the dependent rotation formulas recompute.
$ws.Range("H2").Value = 90
$ws.Range("L2").Value = 90
$ws.Range("H5").Value = 225

# petit commentaire / effacement commentaire: scroll the view over to
# column F and move the selection (previously K8) to I11.
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 6
$aw.ScrollRow = 1
$ws.Range("I11").Select()
